# The 2009年 data row (row 2) is removed from the sheet; subsequent rows
# shift up by one (2010年 -> row 2, 2011年 -> row 3), and the sheet's
# used range shrinks from A1:P4 to A1:P3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("2:2").Delete()
